$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values regenerated for s_vals (filtered save games), per row (B,C,D,E,G) - F is unchanged
$data = @{
    2  = @{ B = 0.7287194209349384; C = 1.65323645889881;  D = 3.082599426703578;  E = 6.48142807727062;  G = 11.94598338380795 }
    3  = @{ B = 0.7287194209349384; C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 3.594575437922795 }
    4  = @{ B = 0.7287194209349384; C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 3.594575437922795 }
    5  = @{ B = 0.1554434735375247; C = 1.65323645889881;  D = 0.7127328510149897; E = 6.48142807727062;  G = 9.002840860721944 }
    6  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    7  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 3.082599426703578;  E = 6.48142807727062;  G = 14.40014219143469 }
    8  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 16.98373111632243; E = 0.4998867070740569; G = 22.31973251085698 }
    9  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 16.98373111632243; E = 0.4998867070740569; G = 22.31973251085698 }
    10 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
